$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "LIVEHTA_723 - Test723" test-case identifiers with the
# new "NewImportLogic_3 - Test_Automation_3" ones (Oncology pops test data
# update). Only the cells that actually contain that identifier change;
# everything else stays the same.

$ws.Range("B2").Value = "NewImportLogic_3 - Test_Automation_3 - 1/13/2023"
$ws.Range("C2").Value = "NewImportLogic_3 - Test_Automation_3"
$ws.Range("D2").Value = "NewImportLogic_3 - Test_Automation_3_radio_button"

$ws.Range("I3").Value = "ExcelReport-NewImportLogic_3 - Test_Automation_3-Clinical-"
$ws.Range("I4").Value = "WordReport-NewImportLogic_3 - Test_Automation_3-Clinical-"
$ws.Range("I5").Value = "ExcelReport-NewImportLogic_3 - Test_Automation_3-Economic-"
$ws.Range("I6").Value = "WordReport-NewImportLogic_3 - Test_Automation_3-Economic-"
$ws.Range("I7").Value = "ExcelReport-NewImportLogic_3 - Test_Automation_3-Quality of Life-"
$ws.Range("I8").Value = "WordReport-NewImportLogic_3 - Test_Automation_3-Quality of Life-"
$ws.Range("I9").Value = "ExcelReport-NewImportLogic_3 - Test_Automation_3-Real-world Evidence-"
$ws.Range("I10").Value = "WordReport-NewImportLogic_3 - Test_Automation_3-Real-world Evidence-"

# Update the view state: scroll so column G is the left-most visible
# column, and move the active selection to I11 (just below the data).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("I11").Select()
